$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row for the "no name" simulado entry.
$ws.Range("A63").Value = "Sem Nome"

# Recreate the built-in "Good" cell style under its Portuguese name ("Bom"),
# matching the author's pt-BR Excel install: light-green fill + dark-green
# font, exactly like the other workbooks in this series.
$bomStyle = $wb.Styles.Add("Bom")
$bomStyle.Font.Color = 24832        # RGB(0,97,0)  -> FF006100
$bomStyle.Interior.Color = 13561798 # RGB(198,239,206) -> FFC6EFCE

$ws.Range("C63").Value = 8.1
$ws.Range("C63").Style = "Bom"
